# "another run of v4" - re-running the Pearson/Spearman correlation script
# overwrites the previously computed correlation coefficients with the
# results of the new run. PearSig (B2) and SpearSig (D2) are unchanged (0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.98022782785691354
$ws.Range("C2").Value = 0.9622607019793904
